$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values could be misinterpreted as numbers by COM
# auto-detection (e.g. "570.22", "0.524"). Force them to Text format first
# so they stay strings, matching the source data, then restore the default
# "Normal" style so no stray number-format style is left behind.
$textCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.870.52"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.134.71"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "570.22"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "150.23"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.134.74"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "0.501"
$ws.Range("E12").Value = "  +6.06%  "
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +9.91%  "
$ws.Range("D14").Value = "37.18"
$ws.Range("E14").Value = "  +5.37%  "
$ws.Range("D15").Value = "3.652.18"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "64.934.44"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "7.17"
$ws.Range("E17").Value = "  +5.71%  "
$ws.Range("D18").Value = "3.134.78"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "509.25"
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  +6.65%  "
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  +7.59%  "
$ws.Range("D23").Value = "15.49"
$ws.Range("E23").Value = "  +11.85%  "
$ws.Range("D24").Value = "7.79"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").Value = "85.03"
$ws.Range("E25").Value = "  +4.58%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "2.91"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("D28").Value = "8.65"
$ws.Range("E28").Value = "  +7.30%  "
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("D30").Value = "27.82"
$ws.Range("E30").Value = "  +5.99%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "2.63"
$ws.Range("E33").Value = "  +5.16%  "
$ws.Range("D34").Value = "6.00"
$ws.Range("E34").Value = "  +7.14%  "
$ws.Range("D35").Value = "6.56"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("D36").Value = "55.57"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "472.91"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").Value = "0.0421"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("D39").Value = "0.0855"
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "3.100.81"
$ws.Range("E41").Value = "  +4.15%  "
$ws.Range("D42").Value = "8.58"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").Value = "0.119"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("E44").Value = "  +10.56%  "
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  +12.19%  "
$ws.Range("D46").Value = "28.96"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").Value = "0.0₃0570"
$ws.Range("E47").Value = "  +10.05%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "0.115"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  +9.41%  "
$ws.Range("D51").Value = "118.67"
$ws.Range("E51").Value = "  -1.92%  "

foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
